$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.439.80'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.05%  '

# Row 16
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.31'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.84%  '

# Row 17
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.115'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.89%  '

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.651.28'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.15%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.418.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.45%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.30%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.73%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.67%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.414.82'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.57%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.486'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.96%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.120'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.88%  '

# Row 11
$ws.Range("E11").Value = '  -9.87%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.374'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.10%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.986.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.91%  '

# Row 14
$ws.Range("E14").Value = '  -9.22%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.573.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.27%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -11.52%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.69%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.88%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '383.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -8.17%  '

# Row 23
$ws.Range("E23").Value = '  -7.54%  '

# Row 24
$ws.Range("E24").Value = '  -0.02%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.61%  '

# Row 26
$ws.Range("E26").Value = '  -2.91%  '

# Row 27
$ws.Range("E27").Value = '  -7.54%  '

# Row 28
$ws.Range("E28").Value = '  +0.34%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.42%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.60%  '

# Row 31
$ws.Range("E31").Value = '  -10.63%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.419.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.69%  '

# Row 33
$ws.Range("E33").Value = '  +0.01%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '22.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.97%  '

# Row 35
$ws.Range("E35").Value = '  -8.81%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '167.80'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.39%  '

# Row 37
$ws.Range("E37").Value = '  -10.54%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.72'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.31%  '

# Row 39
$ws.Range("E39").Value = '  -7.07%  '

# Row 40
$ws.Range("E40").Value = '  -11.30%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0753'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.74%  '

# Row 42
$ws.Range("E42").Value = '  -5.22%  '

# Row 43
$ws.Range("E43").Value = '  -0.52%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.87'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.28%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -13.86%  '

# Row 46
$ws.Range("E46").Value = '  -9.05%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.12'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.93%  '

# Row 48
$ws.Range("E48").Value = '  -2.17%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.43'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.23%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -13.44%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.160.07'
$ws.Range("D51").Style = "Normal"
